# --- Add "2022-Q1" worksheet after "2021-Q4", before "总计" ---
$wb = $excel.ActiveWorkbook
$q4 = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $q4)
$newSheet.Name = "2022-Q1"

# Header row (row 1, columns B:H)
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $newSheet.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Fund holding rows: code, name, fund-scale, total stock position, position ratio, held value (100M), position rank
$data = @(
    @('000118', '广发聚鑫债券A', '208.72', '20.02', '1.02', '2.1289', 9),
    @('213008', '宝盈资源优选混合', '11.31', '81.52', '5.73', '0.6481', 3),
    @('519171', '浦银安盛医疗健康灵活配置混合', '15.50', '88.21', '2.66', '0.4123', 10),
    @('000339', '长城医疗保健混合', '11.06', '88.40', '3.41', '0.3771', 6),
    @('011673', '长城医药科技六个月持有期混合型证券投资基金A', '8.65', '88.50', '3.61', '0.3123', 6),
    @('000780', '鹏华医疗保健股票', '7.67', '82.80', '3.89', '0.2984', 4),
    @('000119', '广发聚鑫债券C', '20.54', '20.02', '1.02', '0.2095', 9),
    @('009623', '长城创新驱动混合', '5.15', '81.88', '3.95', '0.2034', 1),
    @('000968', '广发中证养老产业指数A', '10.39', '94.08', '1.81', '0.1881', 1),
    @('008786', '长城健康生活灵活配置混合', '6.36', '78.57', '2.40', '0.1526', 4),
    @('013037', '长城大健康混合A', '7.41', '39.47', '1.58', '0.1171', 8),
    @('007518', '东方阿尔法优选混合A', '2.03', '72.64', '1.72', '0.0349', 5),
    @('164401', '前海开源中证健康产业指数', '2.13', '94.15', '1.31', '0.0279', 4),
    @('006165', '建信中证1000指数增强A', '2.75', '93.00', '0.99', '0.0272', 6),
    @('011674', '长城医药科技六个月持有期混合型证券投资基金C', '0.72', '88.50', '3.61', '0.0260', 6),
    @('014416', '泰康研究精选股票A', '1.02', '79.61', '2.31', '0.0236', 4),
    @('014417', '泰康研究精选股票C', '0.88', '79.61', '2.31', '0.0203', 4),
    @('002982', '广发中证养老产业指数C', '0.88', '94.08', '1.81', '0.0159', 1),
    @('007519', '东方阿尔法优选混合C', '0.82', '72.64', '1.72', '0.0141', 5),
    @('516560', '华宝养老ETF', '0.75', '97.92', '1.87', '0.0140', 1),
    @('006166', '建信中证1000指数增强C', '0.65', '93.00', '0.99', '0.0064', 6),
    @('013038', '长城大健康混合C', '0.39', '39.47', '1.58', '0.0062', 8),
    @('004641', '万家量化睿选灵活配置混合', '0.16', '85.90', '1.52', '0.0024', 2),
    @('013442', '建信中证1000指数增强E', '0.02', '93.00', '0.99', '0.0002', 6)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $r + 2
    $rec = $data[$r]

    $idxCell = $newSheet.Cells.Item($row, 1)
    $idxCell.Value = $r
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    for ($c = 0; $c -lt 6; $c++) {
        $cell = $newSheet.Cells.Item($row, $c + 2)
        $cell.NumberFormat = '@'
        $cell.Value = $rec[$c]
    }

    $rankCell = $newSheet.Cells.Item($row, 8)
    $rankCell.Value = $rec[6]
}

# --- Update "总计" worksheet: insert a new "2022-Q1" summary row at the top ---
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert(-4121)
$total.Rows.Item(2).ClearFormats()

$totalCell = $total.Cells.Item(2, 1)
$totalCell.Value = 0
$totalCell.Font.Bold = $true
$totalCell.HorizontalAlignment = -4108
$totalCell.VerticalAlignment = -4160
$totalCell.Borders.LineStyle = 1

$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 24
$total.Cells.Item(2, 4).Value = 5.27

# Renumber the index column (A) for the remaining (shifted) rows sequentially
for ($row = 3; $row -le 6; $row++) {
    $total.Cells.Item($row, 1).Value = $row - 2
}

# Restore original active sheet/selection (adding a sheet shifts focus to it)
$wb.Worksheets.Item("2021-Q1").Activate()

